$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 55.81
$ws.Range("B2").Value = 22.12
$ws.Range("B3").Value = 33.23
$ws.Range("B4").Value = 45.99
$ws.Range("B5").Value = 90.01

$ws.Range("C2").Select()
